$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("input_config_2")
$ws.Activate()

$ws.Range("A2").Value = 31868
$ws.Range("B2").Value = "AP"

$ws.Range("A2").Select()
